$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 3677.9546
$ws.Range("I88").Value = 500
$ws.Range("K88").Value = 500
$ws.Range("M88").Value = -94

$ws.Range("H91").Value = 3677.9546
$ws.Range("I91").Value = 500
$ws.Range("K91").Value = 500
$ws.Range("M91").Value = 904

$ws.Range("H113").Value = 4487.4707
$ws.Range("I113").Value = 3154.3333
$ws.Range("J113").Value = 5987.25
$ws.Range("K113").Value = 3154.3333
$ws.Range("L113").Value = 5987.25
$ws.Range("M113").Value = 99.66670000000022
$ws.Range("N113").Value = -12495.25

$ws.Range("H116").Value = 4998.625
$ws.Range("I116").Value = 4150
$ws.Range("J116").Value = 5847.25
$ws.Range("K116").Value = 4150
$ws.Range("L116").Value = 5847.25
$ws.Range("M116").Value = -708
$ws.Range("N116").Value = -12731.25

$ws.Range("H131").Value = 4708.4736
$ws.Range("I131").Value = 2860.0908
$ws.Range("J131").Value = 7250
$ws.Range("K131").Value = 8580.2724
$ws.Range("L131").Value = 21750
$ws.Range("M131").Value = -3540.2724
$ws.Range("N131").Value = -31830

$ws.Range("H138").Value = 3052.8408
$ws.Range("I138").Value = 1176.6086
$ws.Range("K138").Value = 3529.8258
$ws.Range("M138").Value = 1610.1742

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5052102.5
$ws.Range("I2").Value = 11111801
$ws.Range("J2").Value = 2353.8333
$ws.Range("K2").Value = 11111801
$ws.Range("L2").Value = 2353.8333
$ws.Range("M2").Value = -11111688
$ws.Range("N2").Value = -2579.8333

$ws.Range("H22").Value = 2188.6667
$ws.Range("I22").Value = 2426.4
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 2426.4
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -2127.4
$ws.Range("N22").Value = -1598

$ws.Range("H45").Value = 4398911.5
$ws.Range("I45").Value = 6155841.5
$ws.Range("K45").Value = 6155841.5
$ws.Range("M45").Value = -6155464.5

$ws.Range("H92").Value = 75617.2
$ws.Range("J92").Value = 75617.2
$ws.Range("L92").Value = 75617.2
$ws.Range("N92").Value = -80609.2

$ws.Range("H116").Value = 5052102.5
$ws.Range("I116").Value = 11111801
$ws.Range("J116").Value = 2353.8333
$ws.Range("K116").Value = 11111801
$ws.Range("L116").Value = 2353.8333
$ws.Range("M116").Value = -11109507
$ws.Range("N116").Value = -6941.8333

$ws.Range("H122").Value = 835287.6
$ws.Range("I122").Value = 582.3333
$ws.Range("J122").Value = 1605784.9
$ws.Range("K122").Value = 1746.9999
$ws.Range("L122").Value = 4817354.699999999
$ws.Range("M122").Value = 703.0001
$ws.Range("N122").Value = -4822254.699999999

$ws.Range("H132").Value = 3896.22
$ws.Range("I132").Value = 5191.593
$ws.Range("J132").Value = 2375.5652
$ws.Range("K132").Value = 15574.779
$ws.Range("L132").Value = 7126.6956
$ws.Range("M132").Value = -13044.779
$ws.Range("N132").Value = -12186.6956

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5052102.5
$ws.Range("I3").Value = 11111801
$ws.Range("J3").Value = 2353.8333
$ws.Range("K3").Value = 11111801
$ws.Range("L3").Value = 2353.8333
$ws.Range("M3").Value = -11111687
$ws.Range("N3").Value = -2581.8333

$ws.Range("H105").Value = 3969961.8
$ws.Range("I105").Value = 4234526
$ws.Range("K105").Value = 4234526
$ws.Range("M105").Value = -4232779

$ws.Range("H107").Value = 7144622
$ws.Range("I107").Value = 7938360
$ws.Range("K107").Value = 7938360
$ws.Range("M107").Value = -7936440

$ws.Range("H137").Value = 49987.5
$ws.Range("J137").Value = 49987.5
$ws.Range("L137").Value = 49987.5
$ws.Range("N137").Value = -60187.5

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").ClearContents()
$ws.Range("N138").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2000
$ws.Range("I22").Value = 2000
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 2000
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -1650

$ws.Range("H31").Value = 3743.35
$ws.Range("J31").Value = 4920.2812
$ws.Range("L31").Value = 4920.2812
$ws.Range("N31").Value = -5510.2812

$ws.Range("H34").Value = 3743.35
$ws.Range("J34").Value = 4920.2812
$ws.Range("L34").Value = 4920.2812
$ws.Range("N34").Value = -5324.2812

$ws.Range("H86").Value = 13952.28
$ws.Range("J86").Value = 15112.692
$ws.Range("L86").Value = 15112.692
$ws.Range("N86").Value = -17358.692

$ws.Range("H89").Value = 13952.28
$ws.Range("J89").Value = 15112.692
$ws.Range("L89").Value = 75563.45999999999
$ws.Range("N89").Value = -86795.45999999999

$ws.Range("H135").Value = 113246.664
$ws.Range("J135").Value = 113246.664
$ws.Range("L135").Value = 113246.664
$ws.Range("N135").Value = -123386.664

$ws.Range("H138").Value = 110000
$ws.Range("J138").Value = 110000
$ws.Range("L138").Value = 110000
$ws.Range("N138").Value = -120280

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 59347.465
$ws.Range("I12").Value = 148196.5
$ws.Range("J12").Value = 114.77778
$ws.Range("K12").Value = 444589.5
$ws.Range("L12").Value = 344.33334
$ws.Range("M12").Value = -444416.5
$ws.Range("N12").Value = -690.33334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3317294.2
$ws.Range("I102").Value = 4445678
$ws.Range("K102").Value = 4445678
$ws.Range("M102").Value = -4444056

$ws.Range("H107").Value = 969.7
$ws.Range("I107").Value = 1157
$ws.Range("K107").Value = 1157
$ws.Range("M107").Value = 763

$ws.Range("H133").Value = 109992.5
$ws.Range("J133").Value = 109992.5
$ws.Range("L133").Value = 109992.5
$ws.Range("N133").Value = -120112.5

$ws.Range("H135").Value = 80000
$ws.Range("J135").Value = 80000
$ws.Range("L135").Value = 80000
$ws.Range("N135").Value = -90140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1237.8605
$ws.Range("I55").Value = 1256.2858
$ws.Range("J55").Value = 1220.2727
$ws.Range("K55").Value = 1256.2858
$ws.Range("L55").Value = 1220.2727
$ws.Range("M55").Value = -1083.2858
$ws.Range("N55").Value = -1566.2727

$ws.Range("H68").Value = 3978.875
$ws.Range("I68").Value = 3475.8572
$ws.Range("K68").Value = 3475.8572
$ws.Range("M68").Value = -2726.8572

$ws.Range("H71").Value = 3978.875
$ws.Range("I71").Value = 3475.8572
$ws.Range("K71").Value = 17379.286
$ws.Range("M71").Value = -13635.286

$ws.Range("H122").Value = 7120.9165
$ws.Range("I122").Value = 5353.4287
$ws.Range("J122").Value = 9595.4
$ws.Range("K122").Value = 16060.2861
$ws.Range("L122").Value = 28786.2
$ws.Range("M122").Value = -13610.2861
$ws.Range("N122").Value = -33686.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

$ws.Range("H122").Value = 3227.9
$ws.Range("I122").Value = 2324.5454
$ws.Range("J122").Value = 4332
$ws.Range("K122").Value = 6973.6362
$ws.Range("L122").Value = 12996
$ws.Range("M122").Value = -4523.6362
$ws.Range("N122").Value = -17896
